# Naukri_Automation_Details.xlsx - update sample login/application data and
# selected cell, and add a cached display text to the e-mail hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlinks -------------------------------------------------------
# Deleting a hyperlink through a Range-scoped Hyperlinks collection clears
# every hyperlink on the sheet in this runtime, so both are recreated.
# Hyperlinks.Add() also temporarily overwrites the cell's text with the
# TextToDisplay value, so the real cell values are (re)applied afterwards.
$ws.Range("A2").Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add(
    $ws.Range("A2"),
    "mailto:demo123@gmail.com",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "siddhesh.vedre12@gmail.com"
) | Out-Null

$ws.Hyperlinks.Add(
    $ws.Range("B2"),
    "mailto:demoPass@1234",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value
) | Out-Null

# --- Row 2 data ---------------------------------------------------------
$ws.Range("A2").Value = "DemoEmail.gamil.com"
$ws.Range("B2").Value = "DemoPass@123"
$ws.Range("D2").Value = "Mumbai, Navi Mumbai, pune"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 1

# --- Selected cell on the sheet view ------------------------------------
$ws.Range("E6").Select() | Out-Null
